# Badge Adding Condition Added
# Update the "latest badge" Course/Assignment name pair on the STAGE sheet
# (row 2, columns M and N) to the newly generated values, mirroring the
# pattern already used by this test fixture (PortfolioCourse<N> / AssignmentName<N>).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("STAGE")

$ws.Cells.Item(2, 13).Value = "PortfolioCourse47023"
$ws.Cells.Item(2, 14).Value = "AssignmentName47023"
